# Update Bug Reporter sub app
#
# The "problem reporting / privacy" sentence used to be stored as a single
# shared string per language (prefix + the localized word for "privacy"
# baked together). This split it into two pieces - a sentence-prefix string
# (trailing space retained) and a standalone "privacy" word string - so the
# UI can style/localize the word independently. That means a brand new
# "PRIVACY" row is inserted right below the existing "PRP_PRIVACY_FRIENDLY"
# row, and everything below it shifts down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 358 is the key "PRP_PRIVACY_FRIENDLY". Trim the trailing "privacy"
# word (and its translations) off of each language column, leaving the
# sentence prefix (with a trailing space where the word used to continue
# the sentence).
$ws.Range("B358").Value = '问题报告程序不会收集您的'
$ws.Range("C358").Value = 'The problem reporting process does not collect your '
$ws.Range("D358").Value = 'Процесс сообщения о проблемах не собирает ваши '
$ws.Range("E358").Value = 'Le processus de signalement de problèmes ne collecte pas vos '
$ws.Range("F358").Value = 'El proceso de informe de problemas no recopila su '

# Insert a brand-new row right after it (row 359) to hold the standalone
# "privacy" word per language; every row from the old 359 onward shifts
# down by one.
$ws.Rows.Item(359).Insert()

$ws.Range("A359").Value = 'PRIVACY'
$ws.Range("B359").Value = '隐私'
$ws.Range("C359").Value = 'privacy'
$ws.Range("D359").Value = 'конфиденциальность'
$ws.Range("E359").Value = 'confidentialité'
$ws.Range("F359").Value = 'privacidad'
